# add some dungeon resource
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix BgImage ids first
$ws.Range("M5").Value = "castle"
$ws.Range("M6").Value = "viliage"
$ws.Range("M7").Value = "tower"

# Then add the dungeon descriptions
$ws.Range("C5").Value = "失落的古城，在城外就可以听到里面发出的各种奇怪的声音。"
$ws.Range("C6").Value = "这是一个废弃的村落，里面到底有什么样的东西呢。"
$ws.Range("C7").Value = "这是一个充满历史的高大建筑，里面蕴藏着许许多多不为人知的秘密"

# Update active selection to match authored workbook state
$ws.Range("C7").Select()
